$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "23.849.26"
$c.ClearFormats()
$ws.Range("E2").Value = "  +1.41%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.663.36"
$c.ClearFormats()
$ws.Range("E3").Value = "  +1.43%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9964"
$c.ClearFormats()
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("E5").Value = "  +0.00%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "303.88"
$c.ClearFormats()
$ws.Range("E6").Value = "  -0.11%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3812"
$c.ClearFormats()
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3635"
$c.ClearFormats()
$ws.Range("E8").Value = "  -0.14%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "51.00"
$c.ClearFormats()
$ws.Range("E9").Value = "  -1.45%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.254"
$c.ClearFormats()
$ws.Range("E10").Value = "  +1.88%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.08220"
$c.ClearFormats()
$ws.Range("E11").Value = "  +0.39%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.9961"
$c.ClearFormats()
$ws.Range("E12").Value = "  -0.19%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "22.77"
$c.ClearFormats()
$ws.Range("E13").Value = "  +1.19%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.545"
$c.ClearFormats()
$ws.Range("E14").Value = "  +1.12%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.474"
$c.ClearFormats()
$ws.Range("E15").Value = "  +1.18%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.00001237"
$c.ClearFormats()
$ws.Range("E16").Value = "  -0.39%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "1.659.93"
$c.ClearFormats()
$ws.Range("E17").Value = "  +1.50%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "97.64"
$c.ClearFormats()
$ws.Range("E18").Value = "  +2.63%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06991"
$c.ClearFormats()
$ws.Range("E19").Value = "  +0.73%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.835"
$c.ClearFormats()
$ws.Range("E20").Value = "  +3.67%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.78"
$c.ClearFormats()
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("E22").Value = "  +0.20%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "12.91"
$c.ClearFormats()
$ws.Range("E23").Value = "  +2.90%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "23.825.23"
$c.ClearFormats()
$ws.Range("E24").Value = "  +1.33%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.511"
$c.ClearFormats()
$ws.Range("E25").Value = "  +0.01%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.076"
$c.ClearFormats()
$ws.Range("E26").Value = "  -0.02%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "21.37"
$c.ClearFormats()
$ws.Range("E27").Value = "  +1.12%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "153.68"
$c.ClearFormats()
$ws.Range("E28").Value = "  +1.23%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.180"
$c.ClearFormats()
$ws.Range("E29").Value = "  -1.81%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "134.52"
$c.ClearFormats()
$ws.Range("E30").Value = "  +0.81%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.843.94"
$c.ClearFormats()
$ws.Range("E31").Value = "  +1.52%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "7.040"
$c.ClearFormats()
$ws.Range("E32").Value = "  +6.11%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.208"
$c.ClearFormats()
$ws.Range("E33").Value = "  +2.04%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.073"
$c.ClearFormats()
$ws.Range("E34").Value = "  +2.12%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "11.79"
$c.ClearFormats()
$ws.Range("E35").Value = "  +3.58%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02834"
$c.ClearFormats()
$ws.Range("E36").Value = "  +2.45%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.2535"
$c.ClearFormats()
$ws.Range("E37").Value = "  +1.62%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "6.139"
$c.ClearFormats()
$ws.Range("E38").Value = "  +1.72%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.08800"
$c.ClearFormats()
$ws.Range("E39").Value = "  +0.24%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.07102"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.50%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "13.15"
$c.ClearFormats()
$ws.Range("E41").Value = "  +7.76%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.7075"
$c.ClearFormats()
$ws.Range("E42").Value = "  +0.65%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.334"
$c.ClearFormats()
$ws.Range("E43").Value = "  -0.66%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "16.16"
$c.ClearFormats()
$ws.Range("E44").Value = "  +1.93%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.6557"
$c.ClearFormats()
$ws.Range("E45").Value = "  +0.53%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.334"
$c.ClearFormats()
$ws.Range("E46").Value = "  +2.25%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.9989"
$c.ClearFormats()
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("E48").Value = "  +0.33%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.07946"
$c.ClearFormats()
$ws.Range("E49").Value = "  -0.37%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "128.71"
$c.ClearFormats()
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.190"
$c.ClearFormats()
$ws.Range("E51").Value = "  -0.39%  "
